# Append a new attendance record as row 25.
# The existing data in the sheet is stored as text (IDs, dates, times are
# all plain strings), so we force a text number format on the cells whose
# content would otherwise be auto-detected as a number/date ("1446896" and
# "2025-01-26"), then restore the "Normal" style once the value has been
# written so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idCell = $ws.Range("A25")
$idCell.NumberFormat = "@"
$idCell.Value = "1446896"
$idCell.Style = "Normal"

$ws.Range("B25").Value = "Asif Newaz"

$dateCell = $ws.Range("C25")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-01-26"
$dateCell.Style = "Normal"

$ws.Range("D25").Value = "23:37:11"
